$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string labels for the added rows
$ws.Range("A6").Value = "user_favcount"
$ws.Range("A7").Value = "user_followercount"
$ws.Range("A8").Value = "user_friendsCount"
$ws.Range("A9").Value = "user_statusCount"

# Row 6 values
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 0.0015
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.3055
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.0032
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 0.00048429
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0.9036

# Row 7 values
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0.1731
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0.1958
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0.2076
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0.26
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0.3649

# Row 8 values
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0.739
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.1582
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0.5821
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0.9223
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0.0919

# Row 9 values
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0.0146
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.1424
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.0257
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 0.0064
$ws.Range("J9").Value = "NaN"
$ws.Range("K9").Value = "NaN"

# Apply the "Good" cell style to the pass columns that are 1 (matching existing pattern)
$ws.Range("B6").Style = "Good"
$ws.Range("F6").Style = "Good"
$ws.Range("H6").Style = "Good"
$ws.Range("B9").Style = "Good"
$ws.Range("F9").Style = "Good"
$ws.Range("H9").Style = "Good"

# I6 uses the scientific-number-format style (same as C4/E4/G4/I4)
$ws.Range("I6").NumberFormat = "0.00E+00"

# Update the selection to match the final state
$ws.Range("H16").Select()
